# Apply the "Update included example data" changes to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "Wheel assembly"

# Update header row
$ws.Range("A1").Value = "PN"
$ws.Range("B1").Value = "QTY"

# Update part numbers in column A
$ws.Range("A2").Value = "SK1006-01"
$ws.Range("A3").Value = "SK1001-01"

# Add a new row of data
$ws.Range("A4").Value = "SK1007-01"
$ws.Range("B4").Value = 2

# Update selection to A5
$ws.Range("A5").Select()
